$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row: new columns AD, AE, AF (values + formatting matching existing header style)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-48: season record values (Wins, Losses, Ties)
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 100  # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 61   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
